$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 55.25
$ws.Range("I5").Value = 37.4
$ws.Range("J5").Value = 85
$ws.Range("K5").Value = 37.4
$ws.Range("L5").Value = 85
$ws.Range("M5").Value = 77.59999999999999
$ws.Range("N5").Value = -315
$ws.Range("H28").Value = 90909480
$ws.Range("I28").Value = 90909480
$ws.Range("K28").Value = 90909480
$ws.Range("M28").Value = -90908995
$ws.Range("H32").Value = 2028.7646
$ws.Range("I32").Value = 1499.909
$ws.Range("J32").Value = 2998.3333
$ws.Range("K32").Value = 1499.909
$ws.Range("L32").Value = 2998.3333
$ws.Range("M32").Value = -1173.909
$ws.Range("N32").Value = -3650.3333
$ws.Range("H41").Value = 991.3333
$ws.Range("I41").Value = 991.3333
$ws.Range("K41").Value = 991.3333
$ws.Range("M41").Value = -551.3333
$ws.Range("H98").Value = 882.04346
$ws.Range("I98").Value = 870.8570999999999
$ws.Range("K98").Value = 870.8570999999999
$ws.Range("M98").Value = 627.1429000000001
$ws.Range("H111").Value = 3561.6667
$ws.Range("I111").Value = 3092.5
$ws.Range("K111").Value = 9277.5
$ws.Range("M111").Value = -6210.5
$ws.Range("H122").Value = 882.04346
$ws.Range("I122").Value = 870.8570999999999
$ws.Range("K122").Value = 2612.5713
$ws.Range("M122").Value = -162.5712999999996
$ws.Range("H132").Value = 920.7895
$ws.Range("I132").Value = 958.82355
$ws.Range("K132").Value = 2876.47065
$ws.Range("M132").Value = -346.4706499999998
$ws.Range("H138").Value = 2270.3333
$ws.Range("I138").Value = 1841.96
$ws.Range("J138").Value = 2735.9565
$ws.Range("K138").Value = 5525.88
$ws.Range("L138").Value = 8207.869499999999
$ws.Range("M138").Value = -385.8800000000001
$ws.Range("N138").Value = -18487.8695
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 37039016
$ws.Range("I2").Value = 55556388
$ws.Range("K2").Value = 55556388
$ws.Range("M2").Value = -55556275
$ws.Range("H32").Value = 5097.635
$ws.Range("I32").Value = 2502.875
$ws.Range("J32").Value = 13400.866
$ws.Range("K32").Value = 2502.875
$ws.Range("L32").Value = 13400.866
$ws.Range("M32").Value = -2215.875
$ws.Range("N32").Value = -13974.866
$ws.Range("H110").Value = 9261610
$ws.Range("I110").Value = 9261610
$ws.Range("K110").Value = 9261610
$ws.Range("M110").Value = -9259565
$ws.Range("H116").Value = 37039016
$ws.Range("I116").Value = 55556388
$ws.Range("K116").Value = 55556388
$ws.Range("M116").Value = -55554094
$ws.Range("H132").Value = 2192.6
$ws.Range("I132").Value = 1992.0714
$ws.Range("K132").Value = 5976.2142
$ws.Range("M132").Value = -3446.2142
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 37039016
$ws.Range("I3").Value = 55556388
$ws.Range("K3").Value = 55556388
$ws.Range("M3").Value = -55556274
$ws.Range("H134").Value = 1653.6842
$ws.Range("I134").Value = 1135.6875
$ws.Range("K134").Value = 3407.0625
$ws.Range("M134").Value = -872.0625
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 2506399
$ws.Range("I6").Value = 632998.75
$ws.Range("J6").Value = 10000000
$ws.Range("K6").Value = 632998.75
$ws.Range("L6").Value = 10000000
$ws.Range("M6").Value = -632885.75
$ws.Range("N6").Value = -10000226
$ws.Range("H31").Value = 2956.125
$ws.Range("I31").Value = 1483.2174
$ws.Range("J31").Value = 6720.222
$ws.Range("K31").Value = 1483.2174
$ws.Range("L31").Value = 6720.222
$ws.Range("M31").Value = -1188.2174
$ws.Range("N31").Value = -7310.222
$ws.Range("H34").Value = 2956.125
$ws.Range("I34").Value = 1483.2174
$ws.Range("J34").Value = 6720.222
$ws.Range("K34").Value = 1483.2174
$ws.Range("L34").Value = 6720.222
$ws.Range("M34").Value = -1281.2174
$ws.Range("N34").Value = -7124.222
$ws.Range("H62").Value = 2561.6667
$ws.Range("I62").Value = 2574
$ws.Range("J62").Value = 2500
$ws.Range("K62").Value = 2574
$ws.Range("L62").Value = 2500
$ws.Range("M62").Value = -1950
$ws.Range("N62").Value = -3748
$ws.Range("H65").Value = 2561.6667
$ws.Range("I65").Value = 2574
$ws.Range("J65").Value = 2500
$ws.Range("K65").Value = 12870
$ws.Range("L65").Value = 12500
$ws.Range("M65").Value = -9750
$ws.Range("N65").Value = -18740
$ws.Range("H69").Value = 24997.5
$ws.Range("I69").Value = 24997.5
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 24997.5
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -24248.5
$ws.Range("N69").Value = $null
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").Value = $null
$ws.Range("H72").Value = 24997.5
$ws.Range("I72").Value = 24997.5
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 74992.5
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -71248.5
$ws.Range("N72").Value = $null
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").Value = $null
$ws.Range("H96").Value = 5054
$ws.Range("J96").Value = 5054
$ws.Range("L96").Value = 5054
$ws.Range("N96").Value = -10546
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 212.64285
$ws.Range("I12").Value = 163.22223
$ws.Range("K12").Value = 489.66669
$ws.Range("M12").Value = -316.66669
$ws.Range("H23").Value = 125064.25
$ws.Range("I23").Value = 68
$ws.Range("J23").Value = 142920.86
$ws.Range("K23").Value = 204
$ws.Range("L23").Value = 428762.58
$ws.Range("M23").Value = 31
$ws.Range("N23").Value = -429232.58
$ws.Range("H38").Value = 57.23077
$ws.Range("J38").Value = 44.555557
$ws.Range("L38").Value = 133.666671
$ws.Range("N38").Value = -827.666671
$ws.Range("H75").Value = 371
$ws.Range("I75").Value = 194.66667
$ws.Range("K75").Value = 584.00001
$ws.Range("M75").Value = 413.99999
$ws.Range("H78").Value = 371
$ws.Range("I78").Value = 194.66667
$ws.Range("K78").Value = 1752.00003
$ws.Range("M78").Value = 3239.99997
$ws.Range("H92").Value = 234.27272
$ws.Range("I92").Value = 226.57143
$ws.Range("K92").Value = 679.71429
$ws.Range("M92").Value = 568.28571
$ws.Range("H107").Value = 56125.555
$ws.Range("I107").Value = 440.66666
$ws.Range("J107").Value = 83968
$ws.Range("K107").Value = 1321.99998
$ws.Range("L107").Value = 251904
$ws.Range("M107").Value = 598.0000199999999
$ws.Range("N107").Value = -255744
$ws.Range("H113").Value = 1199.6666
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1399.5
$ws.Range("K113").Value = 2400
$ws.Range("L113").Value = 4198.5
$ws.Range("M113").Value = -230
$ws.Range("N113").Value = -8538.5
$ws.Range("H121").Value = 546.5714
$ws.Range("I121").Value = 373.25
$ws.Range("J121").Value = 777.6667
$ws.Range("K121").Value = 1119.75
$ws.Range("L121").Value = 2333.0001
$ws.Range("M121").Value = 190.25
$ws.Range("N121").Value = -4953.0001
$ws.Range("H131").Value = 1365.9375
$ws.Range("I131").Value = 497.5
$ws.Range("K131").Value = 1492.5
$ws.Range("M131").Value = 3547.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H17").Value = 8332.666999999999
$ws.Range("I17").Value = 4999
$ws.Range("K17").Value = 4999
$ws.Range("M17").Value = -4829
$ws.Range("H18").Value = 29999
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null
$ws.Range("H21").Value = 7999.625
$ws.Range("J21").Value = 7999.625
$ws.Range("L21").Value = 7999.625
$ws.Range("N21").Value = -8347.625
$ws.Range("H22").Value = 3358.5715
$ws.Range("J22").Value = 2227.1
$ws.Range("L22").Value = 2227.1
$ws.Range("N22").Value = -2817.1
$ws.Range("H23").Value = 7011335
$ws.Range("I23").Value = 7011335
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 7011335
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = -7011105
$ws.Range("N23").Value = $null
$ws.Range("H24").Value = 9999.333000000001
$ws.Range("J24").Value = 9999.5
$ws.Range("L24").Value = 9999.5
$ws.Range("N24").Value = -10685.5
$ws.Range("H27").Value = 3358.5715
$ws.Range("J27").Value = 2227.1
$ws.Range("L27").Value = 2227.1
$ws.Range("N27").Value = -2441.1
$ws.Range("H30").Value = 6500
$ws.Range("J30").Value = 7000
$ws.Range("L30").Value = 7000
$ws.Range("N30").Value = -7216
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").Value = $null
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").Value = $null
$ws.Range("H122").Value = 3706.3635
$ws.Range("I122").Value = 3702.7
$ws.Range("J122").Value = 3743
$ws.Range("K122").Value = 11108.1
$ws.Range("L122").Value = 11229
$ws.Range("M122").Value = -8658.099999999999
$ws.Range("N122").Value = -16129
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = $null
$ws.Range("H113").Value = 607.875
$ws.Range("J113").Value = 593
$ws.Range("L113").Value = 1779
$ws.Range("N113").Value = -6119
